$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Added Indian MF 1st Stab" update.
#
# The report adds 9 new weekly date columns (Jun_16 .. Sep_08) in front of
# the existing data: they land in columns B..J, and everything that used to
# start at column B now reads starting at column K for the header row and
# the one data row that carries a rating-change note (row 5, Zacks
# Investment Research). All of the other analyst rows simply gain 9 more
# "UN" cells tacked on after their current last column (their existing
# B..V content is left exactly where it is).
# ---------------------------------------------------------------------------

# ----- Row 1: header dates, newest (col B) to oldest (col AE) --------------
$row1Vals = @("Sep_08","Aug_25","Aug_04","Jul_23","Jul_17","Jul_07","Jun_30","Jun_24","Jun_16","Jun_09","Jun_03","May_27","May_23","May_19","May_15","May_12","May_05","Apr_28","Apr_24","Apr_21","Apr_17","Apr_11","Apr_07","Apr_04","Mar_31","Mar_27","Mar_24","Mar_17","Mar_13","Mar_10")
for ($i = 0; $i -lt $row1Vals.Length; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $row1Vals[$i]
}

# ----- Row 5: Zacks Investment Research ------------------------------------
# A5 (firm name) is untouched. B5..AE5 default to "UN" except:
#   F5  -> new 7/17/2019 Downgrade note (Strong-Buy -> Hold), pink highlight
#   X5  -> the pre-existing 4/5/2019 Upgrade note (Hold -> Buy, $3.25), which
#          used to live at O5 before the 9 new columns pushed it over,
#          keeping its original green highlight
for ($c = 2; $c -le 31; $c++) {
    $ws.Cells.Item(5, $c).Value = "UN"
}
$noteCell = $ws.Cells.Item(5, 6)
$noteCell.Value = "7/17/2019,Downgrades,Strong-Buy -> Hold,"
$noteCell.Interior.Color = 13353215

$oldNoteCell = $ws.Cells.Item(5, 24)
$oldNoteCell.Value = "4/5/2019,Upgrades,Hold -> Buy,`$3.25"
$oldNoteCell.Interior.Color = 13434828

# ----- All other analyst rows: append 9 more "UN" cells ---------------------
# Each of these rows keeps its existing content untouched; we just extend it
# by 9 cells (matching the 9 new date columns) filled with "UN".
$appendRows = @(2,3,4,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)
foreach ($r in $appendRows) {
    for ($c = 23; $c -le 31; $c++) {
        $ws.Cells.Item($r, $c).Value = "UN"
    }
}

$appendRowsShort1 = @(30,31)
foreach ($r in $appendRowsShort1) {
    for ($c = 17; $c -le 25; $c++) {
        $ws.Cells.Item($r, $c).Value = "UN"
    }
}

$appendRowsShort2 = @(32,33)
foreach ($r in $appendRowsShort2) {
    for ($c = 8; $c -le 16; $c++) {
        $ws.Cells.Item($r, $c).Value = "UN"
    }
}

# ----- Column widths: keep the uniform 8.0 width across the newly used
# columns (matches the existing C..V formatting being carried through to
# the new columns up to AE). -------------------------------------------------
for ($c = 23; $c -le 31; $c++) {
    $ws.Columns($c).ColumnWidth = 7.14
}
